# Refresh the cryptos table on Sheet1 (GitHub Actions data update).
# Column D ("Price") holds values as literal text in the source file. Plain
# decimal-looking strings (e.g. "316.77") would otherwise be auto-converted to
# numbers by Excel when assigned to a General-formatted cell, so those get a
# leading apostrophe to force text entry, exactly like typing them in by hand.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.500.11'
$ws.Range("E2").Value = '  -0.54%  '

$ws.Range("D3").Value = '1.879.70'
$ws.Range("E3").Value = '  -1.32%  '

$ws.Range("D4").Value = '''1.016'
$ws.Range("E4").Value = '  -1.02%  '

$ws.Range("D5").Value = '''316.77'
$ws.Range("E5").Value = '  -0.97%  '

$ws.Range("D6").Value = '''1.014'
$ws.Range("E6").Value = '  -1.25%  '

$ws.Range("D7").Value = '''0.5114'
$ws.Range("E7").Value = '  -1.67%  '

$ws.Range("D8").Value = '''0.3950'
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").Value = '''0.08420'
$ws.Range("E9").Value = '  +0.82%  '

$ws.Range("D10").Value = '''1.111'
$ws.Range("E10").Value = '  -2.07%  '

$ws.Range("D11").Value = '''6.272'
$ws.Range("E11").Value = '  -0.48%  '

$ws.Range("D12").Value = '1.886.54'
$ws.Range("E12").Value = '  -1.61%  '

$ws.Range("D13").Value = '''20.49'
$ws.Range("E13").Value = '  -0.85%  '

$ws.Range("D14").Value = '''7.285'
$ws.Range("E14").Value = '  -0.38%  '

$ws.Range("D15").Value = '''1.016'
$ws.Range("E15").Value = '  -1.28%  '

$ws.Range("D16").Value = '''0.00001109'
$ws.Range("E16").Value = '  -0.45%  '

$ws.Range("D17").Value = '''91.29'
$ws.Range("E17").Value = '  -0.53%  '

$ws.Range("D18").Value = '''0.06760'
$ws.Range("E18").Value = '  -0.66%  '

$ws.Range("D19").Value = '''17.72'
$ws.Range("E19").Value = '  -1.49%  '

$ws.Range("D20").Value = '''1.013'
$ws.Range("E20").Value = '  -1.37%  '

$ws.Range("D21").Value = '''5.961'
$ws.Range("E21").Value = '  -2.08%  '

$ws.Range("D22").Value = '28.517.97'
$ws.Range("E22").Value = '  -0.71%  '

$ws.Range("E23").Value = '  -1.15%  '

$ws.Range("D24").Value = '''2.274'
$ws.Range("E24").Value = '  -0.66%  '

$ws.Range("D25").Value = '2.101.82'
$ws.Range("E25").Value = '  -1.52%  '

$ws.Range("D26").Value = '''161.28'
$ws.Range("E26").Value = '  -0.80%  '

$ws.Range("D27").Value = '''20.86'
$ws.Range("E27").Value = '  -0.65%  '

$ws.Range("E28").Value = '  -2.86%  '

$ws.Range("D29").Value = '''127.55'
$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("E30").Value = '  -0.47%  '

$ws.Range("D31").Value = '''1.049'
$ws.Range("E31").Value = '  -0.54%  '

$ws.Range("D32").Value = '''5.790'
$ws.Range("E32").Value = '  -2.80%  '

$ws.Range("D33").Value = '''3.618'
$ws.Range("E33").Value = '  -1.66%  '

$ws.Range("D34").Value = '''0.02440'
$ws.Range("E34").Value = '  -1.14%  '

$ws.Range("E35").Value = '  -2.34%  '

$ws.Range("D36").Value = '''0.2185'
$ws.Range("E36").Value = '  -1.57%  '

$ws.Range("D37").Value = '''8.951'
$ws.Range("E37").Value = '  -5.09%  '

$ws.Range("D38").Value = '''1.267'
$ws.Range("E38").Value = '  +0.60%  '

$ws.Range("D39").Value = '''1.193'
$ws.Range("E39").Value = '  -0.19%  '

$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = '''5.073'
$ws.Range("E40").Value = '  +1.03%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '''0.6452'
$ws.Range("E41").Value = '  -1.74%  '

$ws.Range("E42").Value = '  +0.33%  '

$ws.Range("E43").Value = '  -1.65%  '

$ws.Range("D44").Value = '''0.6075'
$ws.Range("E44").Value = '  -1.46%  '

$ws.Range("D45").Value = '''13.09'
$ws.Range("E45").Value = '  -1.79%  '

$ws.Range("D46").Value = '''3.716'
$ws.Range("E46").Value = '  -0.90%  '

$ws.Range("D47").Value = '''2.018'
$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("D48").Value = '''1.206'
$ws.Range("E48").Value = '  -7.39%  '

$ws.Range("E49").Value = '  -2.19%  '

$ws.Range("D50").Value = '''122.41'
$ws.Range("E50").Value = '  -0.40%  '

$ws.Range("D51").Value = '''0.06849'
$ws.Range("E51").Value = '  -1.47%  '
